# Updates the "cryptos" sheet with refreshed price/volume data (and swaps the
# Toncoin/Monero rows), matching the GitHub Actions data-refresh commit.
#
# Columns D (Price) and E (Volume(1h)) are plain text cells in the workbook
# (t="inlineStr"), not numbers/percentages. Many of the new "Price" strings
# look like plain decimals (e.g. "1.013"), and Excel's COM layer will happily
# auto-convert a bare numeric-looking string assigned via .Value into a real
# number, which would corrupt the cell type. To avoid that, those values are
# assigned with a leading apostrophe (the standard Excel "treat as text"
# quote-prefix), exactly as a user typing '1.013 into the cell would do.
# Values that are not valid numbers (e.g. "27.881.02", containing two dots)
# are assigned directly since Excel already treats them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.881.02'
$ws.Range("E2").Value = '  +2.77%  '

$ws.Range("D3").Value = '1.870.55'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").Value = '''1.013'
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").Value = '''313.60'
$ws.Range("E5").Value = '  +1.11%  '

$ws.Range("D6").Value = '''1.012'
$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("D7").Value = '''0.4828'
$ws.Range("E7").Value = '  +1.02%  '

$ws.Range("D8").Value = '''0.3823'
$ws.Range("E8").Value = '  +3.45%  '

$ws.Range("D9").Value = '''0.07381'
$ws.Range("E9").Value = '  +1.54%  '

$ws.Range("D10").Value = '''0.9399'
$ws.Range("E10").Value = '  +0.74%  '

$ws.Range("D11").Value = '''21.03'
$ws.Range("E11").Value = '  +5.53%  '

$ws.Range("D12").Value = '''0.07819'

$ws.Range("D13").Value = '1.888.36'
$ws.Range("E13").Value = '  +1.85%  '

$ws.Range("D14").Value = '''5.496'
$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").Value = '''6.609'
$ws.Range("E15").Value = '  +1.76%  '

$ws.Range("D16").Value = '''91.02'
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("E17").Value = '  -0.41%  '

$ws.Range("D18").Value = '''0.000008870'
$ws.Range("E18").Value = '  +2.02%  '

$ws.Range("E19").Value = '  -0.46%  '

$ws.Range("D20").Value = '27.888.65'
$ws.Range("E20").Value = '  +2.73%  '

$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").Value = '2.126.73'
$ws.Range("E23").Value = '  +2.79%  '

$ws.Range("D24").Value = '''10.84'
$ws.Range("E24").Value = '  +1.70%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''156.83'
$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '''1.935'
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").Value = '''18.57'
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("D28").Value = '''2.057'
$ws.Range("E28").Value = '  +3.46%  '

$ws.Range("D29").Value = '''116.08'
$ws.Range("E29").Value = '  +0.99%  '

$ws.Range("D30").Value = '''4.984'
$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("D31").Value = '''0.08920'
$ws.Range("E31").Value = '  +0.45%  '

$ws.Range("D32").Value = '''3.336'
$ws.Range("E32").Value = '  +0.94%  '

$ws.Range("D33").Value = '''1.221'
$ws.Range("E33").Value = '  +3.44%  '

$ws.Range("D34").Value = '''0.7654'
$ws.Range("E34").Value = '  +3.78%  '

$ws.Range("D35").Value = '''4.652'
$ws.Range("E35").Value = '  +2.68%  '

$ws.Range("D36").Value = '''2.716'
$ws.Range("E36").Value = '  +1.14%  '

$ws.Range("D37").Value = '''1.136'
$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("E38").Value = '  +3.00%  '

$ws.Range("D39").Value = '''0.5658'
$ws.Range("E39").Value = '  +6.97%  '

$ws.Range("D40").Value = '''0.05372'
$ws.Range("E40").Value = '  +2.01%  '

$ws.Range("D41").Value = '''2.999'
$ws.Range("E41").Value = '  +0.63%  '

$ws.Range("D42").Value = '''7.057'
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("D43").Value = '''8.587'
$ws.Range("E43").Value = '  +3.33%  '

$ws.Range("E44").Value = '  +0.58%  '

$ws.Range("D45").Value = '''0.4896'
$ws.Range("E45").Value = '  +3.15%  '

$ws.Range("D46").Value = '''10.71'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").Value = '''105.36'
$ws.Range("E47").Value = '  +3.20%  '

$ws.Range("D48").Value = '''1.013'
$ws.Range("E48").Value = '  -0.42%  '

$ws.Range("D49").Value = '''1.674'
$ws.Range("E49").Value = '  +3.12%  '

$ws.Range("D50").Value = '''67.68'
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("D51").Value = '''0.06110'
$ws.Range("E51").Value = '  +0.83%  '
